$d = $word.ActiveDocument
$CR = [char]13

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replace failed for: $old"
    }
}

# 1. Fix typo "assingment" -> "assignment"
$old = "In my solution to this assingment, I have included 8 feature sets. They are listed below:"
$new = "In my solution to this assignment, I have included 8 feature sets. They are listed below:"
Replace-Text $old $new

# 2. Replace the "Besides the MaxEnt tagging..." paragraph with the new "Please note..." sentence.
$old = "Besides the MaxEnt tagging, I also have implemented 2 versions of Viterbi and applied them with model trained with feature set 1 and model trained with features set 8, aka. model1 and model8 respectively. The simple version means that the posterior probablity is only calculated using the give feature set, and the complex version means that the posterior probability is calculated with the given feature set augmented with each state (i.e. I, B or O). The complex one should strictly mimic the MEMM in the textbook, and the simple one has some relaxation. However, the simple one shows identical performance of MaxEnt, but the complex one has way too worse performance. Please see results below."
$new = "Please note that during tagging, the nextTag is always set to null to mimic the fact that only prior state has been seen."
Replace-Text $old $new

# 3. Insert one new empty paragraph after it, then a new paragraph with the Viterbi description.
$p = $d.Paragraphs(30)
$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(31)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(32)
$p3.Range.Text = "Besides the MaxEnt tagging, I also have implemented 2 versions of Viterbi and applied them with model trained with features set 8, aka. model8. The complex one should strictly mimic the MEMM in the textbook, and the simple one has some relaxation. The simple one is not 100% correct viterbi implementation (the previous state is always set to null) so it is more or less a placeholder simply for performance comparison purpose. However its performance is not bad at all."

# 4. Update the "Model#: Correct tags: N" lines (text + a couple of numbers changed).
Replace-Text "Model1: Correct tags: 9186" "Model1: correct tags 9186"
Replace-Text "Model2: Correct tags: 9186" "Model2: correct tags 9186"
Replace-Text "Model3: Correct tags: 9356" "Model3: correct tags 9356"
Replace-Text "Model4: Correct tags: 9355" "Model4: correct tags 9355"
Replace-Text "Model5: Correct tags: 9355" "Model5: correct tags 9355"
Replace-Text "Model6: Correct tags: 9351" "Model6: correct tags 9351"
Replace-Text "Model7: Correct tags: 9358" "Model7: correct tags 9345"
Replace-Text "Model8: Correct tags: 9363" "Model8: correct tags 9325"

# 5. Remove the "Model1 with Viterbi-simple/complex" lines entirely.
$old = "Model1 with Viterbi-simple: Correct tags: 9186" + $CR
Replace-Text $old ""
$old = "Model1 with Viterbi-complex: Correct tags: 4480" + $CR
Replace-Text $old ""

# 6. Update the remaining "Model8 with Viterbi-*" lines.
Replace-Text "Model8 with Viterbi-simple: Correct tags: 9363" "Model8 with Viterbi-simple: correct tags 9319"
Replace-Text "Model8 with Viterbi-complex: Correct tags: 4497" "Model8 with Viterbi-complex: correct tags 9331"

# 7. Merge the two "Measures of model5" paragraphs into one.
$old = "Measures of model5: precision:0.8946957878315133 recall:0.9013752455795678 " + $CR + "F1:0.898023096496379"
$new = "Measures of model5: precision:0.8946957878315133 recall:0.9013752455795678 F1:0.898023096496379"
Replace-Text $old $new

# 8. Remove the blank paragraph between "Measures of model6" and "Measures of model7", and update model7's numbers.
$old = "Measures of model6: precision:0.8924814959096221 recall:0.900196463654224 F1:0.8963223787167449 F0.5:0.8940138921407943 F2:0.8986428179179414" + $CR + $CR + "Measures of model7: precision:0.8966731898238748 recall:0.900196463654224 F1:0.8984313725490197 F0.5:0.8973756365060712 F2:0.8994895956026698"
$new = "Measures of model6: precision:0.8924814959096221 recall:0.900196463654224 F1:0.8963223787167449 F0.5:0.8940138921407943 F2:0.8986428179179414" + $CR + "Measures of model7: precision:0.8895729126587149 recall:0.9084479371316306 F1:0.8989113530326595 F0.5:0.8932849084305696 F2:0.904609124344628"
Replace-Text $old $new

# 9. Update model8's measures.
Replace-Text "Measures of model8: precision:0.8997247345654739 recall:0.8990176817288802 F1:0.89937106918239 F0.5:0.8995832350397105 F2:0.8991590033797061" `
             "Measures of model8: precision:0.8867996930161166 recall:0.9080550098231827 F1:0.8973014948553678 F0.5:0.8909707764669595 F2:0.9037228218363835"

# 10. Remove the "Measures of model1 with Viterbi-simple/complex" paragraphs (and their blank separators).
$old = "Measures of model8: precision:0.8867996930161166 recall:0.9080550098231827 F1:0.8973014948553678 F0.5:0.8909707764669595 F2:0.9037228218363835" + $CR + $CR + "Measures of model1 with Viterbi-simple: precision:0.8182861514919664 recall:0.8404715127701375 F1:0.8292304710215158 F0.5:0.8226290285362664 F2:0.8359387212755979" + $CR + $CR + "Measures of model1 with Viterbi-complex: precision:0.004557538928978352 recall:0.004715127701375246 F1:0.004634994206257241 F0.5:0.004588208304657031 F2:0.004682744088035589" + $CR + $CR + "Measures of model8 with Viterbi-simple"
$new = "Measures of model8: precision:0.8867996930161166 recall:0.9080550098231827 F1:0.8973014948553678 F0.5:0.8909707764669595 F2:0.9037228218363835" + $CR + $CR + "Measures of model8 with Viterbi-simple"
Replace-Text $old $new

# 11. Update the "Measures of model8 with Viterbi-*" lines.
Replace-Text "Measures of model8 with Viterbi-simple: precision:0.8997247345654739 recall:0.8990176817288802 F1:0.89937106918239 F0.5:0.8995832350397105 F2:0.8991590033797061" `
             "Measures of model8 with Viterbi-simple: precision:0.880061115355233 recall:0.9053045186640472 F1:0.8925043579314351 F0.5:0.8849965429822539 F2:0.90014064697609"
Replace-Text "Measures of model8 with Viterbi-complex: precision:0.005050505050505051 recall:0.005108055009823182 F1:0.005079117015042001 F0.5:0.0050619110661163464 F2:0.0050964403324447235" `
             "Measures of model8 with Viterbi-complex: precision:0.8875239923224568 recall:0.9084479371316306 F1:0.8978640776699028 F0.5:0.891631315079059 F2:0.9041845913179508"

# 12. Add an extra blank paragraph right before "Conclusion".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Conclusion")) {
        $d.Paragraphs($i - 1).Range.InsertParagraphAfter()
        break
    }
}

# 13. Update the Conclusion text.
$old = ": Feature conjunction does help improved the performance. Including prior and following states (I assume states here mean BIO tags) as features also helped improve the performance. The best performance I could get is with feature set 8 using MaxEnt and Viterbi-simple using feature set 8. And the strictly mimic Viterbi (i.e.Viterbi-complex) gives surprising low performance. Not sure if I have done anything wrong or not, but I really think I have done it the right way."
$new = ": Feature conjunction does help improved the performance. Including prior state as a feature also helped improve the performance. The best performance I could get is with feature set 8 using MaxEnt and Viterbi-complex using feature set 8. Viterbi did improved the performance but just a little bit."
Replace-Text $old $new

Write-Output "done"
